# "Update Work Week and Social Spending"
#
# The "Data" sheet holds one GDP-per-Capita reading per year for Iraq
# (country code 368). This update refreshes the historical series with
# newer figures and appends six new years (2011-2016).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- 1. Refresh existing year values (1820-2010) with updated figures ---
# Values are re-entered with a leading apostrophe so they are stored as
# text (matching the existing column, which stores its numbers as text)
# instead of being auto-converted to numeric cells.
$changes = @(
    @{Row=2; Value="877"},
    @{Row=52; Value="956"},
    @{Row=95; Value="1275"},
    @{Row=132; Value="2174"},
    @{Row=133; Value="2303"},
    @{Row=134; Value="2482"},
    @{Row=135; Value="3394"},
    @{Row=136; Value="3926"},
    @{Row=137; Value="3663"},
    @{Row=138; Value="3808"},
    @{Row=139; Value="3666"},
    @{Row=140; Value="3974"},
    @{Row=141; Value="4022"},
    @{Row=142; Value="4360"},
    @{Row=143; Value="4720"},
    @{Row=144; Value="4809"},
    @{Row=145; Value="4578"},
    @{Row=146; Value="4965"},
    @{Row=147; Value="5241"},
    @{Row=148; Value="5338"},
    @{Row=149; Value="5043"},
    @{Row=150; Value="5745"},
    @{Row=151; Value="5745"},
    @{Row=152; Value="5536"},
    @{Row=153; Value="5686"},
    @{Row=154; Value="5297"},
    @{Row=155; Value="5982"},
    @{Row=156; Value="6097"},
    @{Row=157; Value="6878"},
    @{Row=158; Value="8007"},
    @{Row=159; Value="7957"},
    @{Row=160; Value="9075"},
    @{Row=161; Value="10769"},
    @{Row=162; Value="10165"},
    @{Row=163; Value="8035"},
    @{Row=164; Value="7704"},
    @{Row=165; Value="6805"},
    @{Row=166; Value="6593"},
    @{Row=167; Value="6268"},
    @{Row=168; Value="5992"},
    @{Row=169; Value="6052"},
    @{Row=170; Value="4635"},
    @{Row=171; Value="4098"},
    @{Row=172; Value="3918"},
    @{Row=173; Value="1518.74208278481"},
    @{Row=174; Value="2046.18088439654"},
    @{Row=175; Value="2132.38232470176"},
    @{Row=176; Value="2108.53800593898"},
    @{Row=177; Value="1996.9277996571"},
    @{Row=178; Value="2337.52265340147"},
    @{Row=179; Value="2337.94787905627"},
    @{Row=180; Value="2794.7864918705"},
    @{Row=181; Value="3157.11246136115"},
    @{Row=182; Value="3649.49294748339"},
    @{Row=183; Value="4147.64851953066"},
    @{Row=184; Value="4631.65253834669"},
    @{Row=185; Value="3782.45263801471"},
    @{Row=186; Value="6061.55461844611"},
    @{Row=187; Value="6441.41021617748"},
    @{Row=188; Value="7117.61350547781"},
    @{Row=189; Value="7587.60567342558"},
    @{Row=190; Value="8587.26471999696"},
    @{Row=191; Value="9268.7091555216"},
    @{Row=192; Value="10274.3303048734"}
)

foreach ($chg in $changes) {
    $ws.Cells.Item($chg.Row, 5).Value = "'" + $chg.Value
}

# --- 2. Append six new years of data (2011-2016) ---
$newRows = @(
    @{Year=2011; Value="11484"},
    @{Year=2012; Value="12652"},
    @{Year=2013; Value="13158"},
    @{Year=2014; Value="12817"},
    @{Year=2015; Value="13014"},
    @{Year=2016; Value="13898"}
)

$startRow = 193
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = 368
    $ws.Cells.Item($r, 2).Value = "Iraq"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $newRows[$i].Year
    $ws.Cells.Item($r, 5).Value = "'" + $newRows[$i].Value
}
